$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Instruction Number values in column B (rows 2-7)
$ws.Range("B2").Value = 3
$ws.Range("B3").Value = 6
$ws.Range("B4").Value = 9
$ws.Range("B5").Value = 8
$ws.Range("B6").Value = 0
$ws.Range("B7").Value = 1

# Update the selected cell/range to B4
$ws.Range("B4").Select()
